$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell holds text data (prices / percentages rendered as strings in the
# source feed). A leading apostrophe forces Excel to keep the literal text
# instead of auto-coercing number-looking values (e.g. "554.79" or "1.00"),
# and resetting the style back to Normal discards the quote-prefix formatting
# flag Excel would otherwise stamp on the cell.
$ws.Range("D2").Value = "'63.442.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.51%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.083.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.24%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'554.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.71%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'137.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.48%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.076.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.53%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.30%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.155"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.08%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.66%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'35.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.07%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.47%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.577.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'63.522.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.30%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.25%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.081.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.33%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'502.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.48%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.00%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.10%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.54%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.14%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'12.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.76%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'76.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.83%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.04%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +1.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'8.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.83%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.37%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.09%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'25.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.01%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.13%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.98%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'526.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -13.44%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'56.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +8.56%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.91%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -6.06%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0410"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.63%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0791"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.43%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.050.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.18%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.43%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'8.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.39%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -10.80%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.251"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.99%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D46").Value = "'2.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.54%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'121.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.47%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'23.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.64%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0₃0491"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -8.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -3.53%  "
$ws.Range("E51").Style = "Normal"
